# Apply updated cryptos list values (price + 1h volume change)
# Source: diff of cryptos.xlsx worksheet XML (rows 2-50, columns D & E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.211.52"
$ws.Range("E2").Value = "  -3.96%  "
$ws.Range("D3").Value = "'2.238.48"
$ws.Range("E3").Value = "  -4.92%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'231.88"
$ws.Range("E5").Value = "  -3.69%  "
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("D7").Value = "'70.33"
$ws.Range("E7").Value = "  -4.58%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.561"
$ws.Range("E9").Value = "  -7.08%  "
$ws.Range("D10").Value = "'0.0999"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'58.28"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "'35.80"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "'6.86"
$ws.Range("E14").Value = "  -6.81%  "
$ws.Range("D15").Value = "'2.571.14"
$ws.Range("E15").Value = "  -4.96%  "
$ws.Range("D16").Value = "'15.04"
$ws.Range("E16").Value = "  -8.54%  "
$ws.Range("E17").Value = "  -5.36%  "
$ws.Range("D18").Value = "'2.243.94"
$ws.Range("E18").Value = "  -4.73%  "
$ws.Range("D19").Value = "'42.084.46"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").Value = "'0.0₃0989"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("D21").Value = "'73.54"
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  -7.62%  "
$ws.Range("D23").Value = "'238.12"
$ws.Range("E23").Value = "  -7.43%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("E27").Value = "  -6.50%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("E29").Value = "  -8.81%  "
$ws.Range("D30").Value = "'168.29"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("E31").Value = "  -9.11%  "
$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "  -7.79%  "
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'4.80"
$ws.Range("E36").Value = "  -8.10%  "
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'22.38"
$ws.Range("E38").Value = "  +17.12%  "
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("D41").Value = "'0.0266"
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("D42").Value = "'66.86"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'4.90"
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("D44").Value = "'8.97"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("E45").Value = "  -9.66%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  -7.28%  "
$ws.Range("D48").Value = "'10.21"
$ws.Range("E48").Value = "  +7.19%  "
$ws.Range("D49").Value = "'4.36"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("E50").Value = "  -6.85%  "
